# Update odds values for the week's games (Jogos_da_Semana_FlashScore)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.62
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65

# Row 4
$ws.Range("G4").Value = 3.7
$ws.Range("I4").Value = 2.1
$ws.Range("R4").Value = 1.5
$ws.Range("AD4").Value = 6
$ws.Range("AI4").Value = 9

# Row 5
$ws.Range("I5").Value = 6.5
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.5
$ws.Range("R5").Value = 1.44
$ws.Range("U5").Value = 2.63
$ws.Range("V5").Value = 1.44
$ws.Range("Y5").Value = 10

# Row 6
$ws.Range("O6").Value = 1.5
$ws.Range("P6").Value = 2.5
$ws.Range("R6").Value = 1.47

# Row 7
$ws.Range("G7").Value = 1.79
$ws.Range("R7").Value = 1.5

# Row 9
$ws.Range("G9").Value = 2.25
$ws.Range("I9").Value = 3.5
$ws.Range("J9").Value = 3.1
$ws.Range("K9").Value = 1.95
$ws.Range("X9").Value = 9.5
$ws.Range("AH9").Value = 8
$ws.Range("AI9").Value = 15
$ws.Range("AW9").Value = 5
$ws.Range("AZ9").Value = 67
$ws.Range("BA9").Value = 101

# Row 12
$ws.Range("G12").Value = 1.73
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 5.5
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 6.5
$ws.Range("W12").Value = 5.5
$ws.Range("X12").Value = 7
$ws.Range("Z12").Value = 13
$ws.Range("AJ12").Value = 19
$ws.Range("AT12").Value = 2.38
$ws.Range("AU12").Value = 9.5
$ws.Range("AW12").Value = 6.5
$ws.Range("BA12").Value = 151

# Row 13
$ws.Range("Q13").Value = 1.93
$ws.Range("R13").Value = 1.93
